# Insert a new weekly price record as row 42, pushing the existing
# rows 42-56 down to 43-57 (dimension grows from A1:R56 to A1:R57).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 42..56 down by one to make room for the new record.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly record.
$ws.Range("A42").Value = 10
$ws.Range("B42").Value = "Vega Modelo de Temuco"
$ws.Range("C42").Value = "La Araucanía"
$ws.Range("D42").Value = 44468
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = 100112035
$ws.Range("G42").Value = "Bruselas (repollito)"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 40
$ws.Range("K42").Value = 25000
$ws.Range("L42").Value = 25000
$ws.Range("M42").Value = 25000
$ws.Range("N42").Value = "$/malla 10 kilos"
$ws.Range("O42").Value = "Provincia de Quillota"
$ws.Range("P42").Value = 2500
$ws.Range("Q42").Value = 10
$ws.Range("R42").Value = "Hortaliza"
